$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StatDef")

# Row 259
$ws.Range("A259").Value = 4256
$ws.Range("B259").Value = "COOKIE"
$ws.Range("C259").Value = "Cookie"
$ws.Range("D259").Value = 25
$ws.Range("E259").Value = 100
$ws.Range("F259").Value = 100
$ws.Range("G259").Value = 100
$ws.Range("H259").Value = 100
$ws.Range("I259").Value = 100
$ws.Range("J259").Value = 100
$ws.Range("K259").Value = 100
$ws.Range("L259").Value = 100
$ws.Range("M259").Value = 10
$ws.Range("N259").Value = 1
$ws.Range("O259").Value = 100
$ws.Range("P259").Value = 100
$ws.Range("Q259").Value = 100
$ws.Range("R259").Value = 100
$ws.Range("S259").Value = 10
$ws.Range("T259").Value = 12
$ws.Range("U259").Value = "Small"
$ws.Range("V259").Value = "Demihuman"
$ws.Range("W259").Value = "Neutral3"
$ws.Range("X259").Value = 1036
$ws.Range("Y259").Value = 240
$ws.Range("Z259").Value = 936
$ws.Range("AA259").Value = 200
$ws.Range("AB259").Value = "Normal"
$ws.Range("AC259").Value = "Normal"
$ws.Range("AD259").Value = "AiPassive"
$ws.Range("AE259").Value = 111
$ws.Range("AF259").Value = "cookie.spr"
$ws.Range("AG259").Value = 0
$ws.Range("AH259").Value = 0.5
$ws.Range("AI259").Value = 1

# Row 260
$ws.Range("A260").Value = 4257
$ws.Range("B260").Value = "COOKIE_XMAS"
$ws.Range("C260").Value = "Cookie"
$ws.Range("D260").Value = 28
$ws.Range("E260").Value = 100
$ws.Range("F260").Value = 100
$ws.Range("G260").Value = 100
$ws.Range("H260").Value = 100
$ws.Range("I260").Value = 100
$ws.Range("J260").Value = 100
$ws.Range("K260").Value = 100
$ws.Range("L260").Value = 100
$ws.Range("M260").Value = 10
$ws.Range("N260").Value = 1
$ws.Range("O260").Value = 100
$ws.Range("P260").Value = 100
$ws.Range("Q260").Value = 100
$ws.Range("R260").Value = 100
$ws.Range("S260").Value = 10
$ws.Range("T260").Value = 12
$ws.Range("U260").Value = "Small"
$ws.Range("V260").Value = "Demihuman"
$ws.Range("W260").Value = "Holy2"
$ws.Range("X260").Value = 1248
$ws.Range("Y260").Value = 240
$ws.Range("Z260").Value = 1248
$ws.Range("AA260").Value = 400
$ws.Range("AB260").Value = "Normal"
$ws.Range("AC260").Value = "Normal"
$ws.Range("AD260").Value = "AiPassive"
$ws.Range("AE260").Value = 350
$ws.Range("AF260").Value = "cookie_xmas.spr"
$ws.Range("AG260").Value = 0
$ws.Range("AH260").Value = 0.5
$ws.Range("AI260").Value = 1

# Row 261
$ws.Range("A261").Value = 4258
$ws.Range("B261").Value = "CRUISER"
$ws.Range("C261").Value = "Cruiser"
$ws.Range("D261").Value = 35
$ws.Range("E261").Value = 100
$ws.Range("F261").Value = 100
$ws.Range("G261").Value = 100
$ws.Range("H261").Value = 100
$ws.Range("I261").Value = 100
$ws.Range("J261").Value = 100
$ws.Range("K261").Value = 100
$ws.Range("L261").Value = 100
$ws.Range("M261").Value = 10
$ws.Range("N261").Value = 7
$ws.Range("O261").Value = 100
$ws.Range("P261").Value = 100
$ws.Range("Q261").Value = 100
$ws.Range("R261").Value = 100
$ws.Range("S261").Value = 10
$ws.Range("T261").Value = 12
$ws.Range("U261").Value = "Medium"
$ws.Range("V261").Value = "Formless"
$ws.Range("W261").Value = "Neutral3"
$ws.Range("X261").Value = 1296
$ws.Range("Y261").Value = 432
$ws.Range("Z261").Value = 1296
$ws.Range("AA261").Value = 400
$ws.Range("AB261").Value = "Normal"
$ws.Range("AC261").Value = "Normal,Ranged"
$ws.Range("AD261").Value = "AiAggressive"
$ws.Range("AE261").Value = 1000
$ws.Range("AF261").Value = "cruiser.spr"
$ws.Range("AG261").Value = 0
$ws.Range("AH261").Value = 0.5
$ws.Range("AI261").Value = 1

# Row 262
$ws.Range("A262").Value = 4259
$ws.Range("B262").Value = "CHEPET"
$ws.Range("C262").Value = "Chepet"
$ws.Range("D262").Value = 42
$ws.Range("E262").Value = 100
$ws.Range("F262").Value = 100
$ws.Range("G262").Value = 100
$ws.Range("H262").Value = 100
$ws.Range("I262").Value = 100
$ws.Range("J262").Value = 100
$ws.Range("K262").Value = 100
$ws.Range("L262").Value = 100
$ws.Range("M262").Value = 10
$ws.Range("N262").Value = 1
$ws.Range("O262").Value = 100
$ws.Range("P262").Value = 100
$ws.Range("Q262").Value = 100
$ws.Range("R262").Value = 100
$ws.Range("S262").Value = 10
$ws.Range("T262").Value = 12
$ws.Range("U262").Value = "Medium"
$ws.Range("V262").Value = "Demihuman"
$ws.Range("W262").Value = "Fire1"
$ws.Range("X262").Value = 672
$ws.Range("Y262").Value = 288
$ws.Range("Z262").Value = 672
$ws.Range("AA262").Value = 400
$ws.Range("AB262").Value = "Normal"
$ws.Range("AC262").Value = "Normal,MiniBoss"
$ws.Range("AD262").Value = "AiAggressive"
$ws.Range("AE262").Value = 350
$ws.Range("AF262").Value = "chepet.spr"
$ws.Range("AG262").Value = 0
$ws.Range("AH262").Value = 0.5
$ws.Range("AI262").Value = 1

# Row 263
$ws.Range("A263").Value = 4260
$ws.Range("B263").Value = "GOBLINE_XMAS"
$ws.Range("C263").Value = "Festive Goblin"
$ws.Range("D263").Value = 25
$ws.Range("E263").Value = 100
$ws.Range("F263").Value = 100
$ws.Range("G263").Value = 100
$ws.Range("H263").Value = 100
$ws.Range("I263").Value = 100
$ws.Range("J263").Value = 100
$ws.Range("K263").Value = 100
$ws.Range("L263").Value = 100
$ws.Range("M263").Value = 10
$ws.Range("N263").Value = 1
$ws.Range("O263").Value = 100
$ws.Range("P263").Value = 100
$ws.Range("Q263").Value = 100
$ws.Range("R263").Value = 100
$ws.Range("S263").Value = 10
$ws.Range("T263").Value = 12
$ws.Range("U263").Value = "Medium"
$ws.Range("V263").Value = "Demihuman"
$ws.Range("W263").Value = "Wind1"
$ws.Range("X263").Value = 1120
$ws.Range("Y263").Value = 240
$ws.Range("Z263").Value = 620
$ws.Range("AA263").Value = 100
$ws.Range("AB263").Value = "Normal"
$ws.Range("AC263").Value = "Normal"
$ws.Range("AD263").Value = "AiAggressive"
$ws.Range("AE263").Value = 400
$ws.Range("AF263").Value = "gobline_xmas.spr"
$ws.Range("AG263").Value = 0
$ws.Range("AH263").Value = 0.5
$ws.Range("AI263").Value = 1

# Row 264
$ws.Range("A264").Value = 4261
$ws.Range("B264").Value = "GARM"
$ws.Range("C264").Value = "Garm"
$ws.Range("D264").Value = 73
$ws.Range("E264").Value = 100
$ws.Range("F264").Value = 100
$ws.Range("G264").Value = 100
$ws.Range("H264").Value = 100
$ws.Range("I264").Value = 100
$ws.Range("J264").Value = 100
$ws.Range("K264").Value = 100
$ws.Range("L264").Value = 100
$ws.Range("M264").Value = 10
$ws.Range("N264").Value = 3
$ws.Range("O264").Value = 100
$ws.Range("P264").Value = 100
$ws.Range("Q264").Value = 100
$ws.Range("R264").Value = 100
$ws.Range("S264").Value = 10
$ws.Range("T264").Value = 12
$ws.Range("U264").Value = "Large"
$ws.Range("V264").Value = "Beast"
$ws.Range("W264").Value = "Water4"
$ws.Range("X264").Value = 608
$ws.Range("Y264").Value = 336
$ws.Range("Z264").Value = 408
$ws.Range("AA264").Value = 400
$ws.Range("AB264").Value = "Boss"
$ws.Range("AC264").Value = "Normal,Elite"
$ws.Range("AD264").Value = "AiAggressive"
$ws.Range("AE264").Value = 325
$ws.Range("AF264").Value = "garm.spr"
$ws.Range("AG264").Value = 0
$ws.Range("AH264").Value = 0.5
$ws.Range("AI264").Value = 1

# Row 265
$ws.Range("A265").Value = 4262
$ws.Range("B265").Value = "GARM_BABY"
$ws.Range("C265").Value = "Baby Garm"
$ws.Range("D265").Value = 61
$ws.Range("E265").Value = 100
$ws.Range("F265").Value = 100
$ws.Range("G265").Value = 100
$ws.Range("H265").Value = 100
$ws.Range("I265").Value = 100
$ws.Range("J265").Value = 100
$ws.Range("K265").Value = 100
$ws.Range("L265").Value = 100
$ws.Range("M265").Value = 10
$ws.Range("N265").Value = 1
$ws.Range("O265").Value = 100
$ws.Range("P265").Value = 100
$ws.Range("Q265").Value = 100
$ws.Range("R265").Value = 100
$ws.Range("S265").Value = 10
$ws.Range("T265").Value = 12
$ws.Range("U265").Value = "Medium"
$ws.Range("V265").Value = "Beast"
$ws.Range("W265").Value = "Water2"
$ws.Range("X265").Value = 879
$ws.Range("Y265").Value = 576
$ws.Range("Z265").Value = 672
$ws.Range("AA265").Value = 450
$ws.Range("AB265").Value = "Normal"
$ws.Range("AC265").Value = "Normal"
$ws.Range("AD265").Value = "AiAggressive"
$ws.Range("AE265").Value = 300
$ws.Range("AF265").Value = "garm_baby.spr"
$ws.Range("AG265").Value = 0
$ws.Range("AH265").Value = 0.5
$ws.Range("AI265").Value = 1

# Row 266
$ws.Range("A266").Value = 4263
$ws.Range("B266").Value = "KNIGHT_OF_WINDSTORM"
$ws.Range("C266").Value = "Stormy Knight"
$ws.Range("D266").Value = 77
$ws.Range("E266").Value = 100
$ws.Range("F266").Value = 100
$ws.Range("G266").Value = 100
$ws.Range("H266").Value = 100
$ws.Range("I266").Value = 100
$ws.Range("J266").Value = 100
$ws.Range("K266").Value = 100
$ws.Range("L266").Value = 100
$ws.Range("M266").Value = 10
$ws.Range("N266").Value = 2
$ws.Range("O266").Value = 100
$ws.Range("P266").Value = 100
$ws.Range("Q266").Value = 100
$ws.Range("R266").Value = 100
$ws.Range("S266").Value = 10
$ws.Range("T266").Value = 12
$ws.Range("U266").Value = "Large"
$ws.Range("V266").Value = "Formless"
$ws.Range("W266").Value = "Wind4"
$ws.Range("X266").Value = 468
$ws.Range("Y266").Value = 288
$ws.Range("Z266").Value = 468
$ws.Range("AA266").Value = 200
$ws.Range("AB266").Value = "Boss"
$ws.Range("AC266").Value = "Normal,WorldBoss"
$ws.Range("AD266").Value = "AiAggressive"
$ws.Range("AE266").Value = 259
$ws.Range("AF266").Value = "knight_of_windstorm.spr"
$ws.Range("AG266").Value = 0
$ws.Range("AH266").Value = 0.5
$ws.Range("AI266").Value = 1

# Row 267
$ws.Range("A267").Value = 4264
$ws.Range("B267").Value = "MYSTCASE"
$ws.Range("C267").Value = "Mystcase"
$ws.Range("D267").Value = 38
$ws.Range("E267").Value = 100
$ws.Range("F267").Value = 100
$ws.Range("G267").Value = 100
$ws.Range("H267").Value = 100
$ws.Range("I267").Value = 100
$ws.Range("J267").Value = 100
$ws.Range("K267").Value = 100
$ws.Range("L267").Value = 100
$ws.Range("M267").Value = 10
$ws.Range("N267").Value = 1
$ws.Range("O267").Value = 100
$ws.Range("P267").Value = 100
$ws.Range("Q267").Value = 100
$ws.Range("R267").Value = 100
$ws.Range("S267").Value = 10
$ws.Range("T267").Value = 12
$ws.Range("U267").Value = "Medium"
$ws.Range("V267").Value = "Formless"
$ws.Range("W267").Value = "Neutral3"
$ws.Range("X267").Value = 1248
$ws.Range("Y267").Value = 432
$ws.Range("Z267").Value = 1248
$ws.Range("AA267").Value = 400
$ws.Range("AB267").Value = "Normal"
$ws.Range("AC267").Value = "Normal"
$ws.Range("AD267").Value = "AiPassive"
$ws.Range("AE267").Value = 450
$ws.Range("AF267").Value = "mystcase.spr"
$ws.Range("AG267").Value = 0
$ws.Range("AH267").Value = 0.5
$ws.Range("AI267").Value = 1

# Row 268
$ws.Range("A268").Value = 4265
$ws.Range("B268").Value = "WRAITH_DEAD"
$ws.Range("AF268").Value = "wraith_dead.spr"
$ws.Range("C268").Value = "Wraith Dead"
$ws.Range("D268").Value = 74
$ws.Range("E268").Value = 100
$ws.Range("F268").Value = 100
$ws.Range("G268").Value = 100
$ws.Range("H268").Value = 100
$ws.Range("I268").Value = 100
$ws.Range("J268").Value = 100
$ws.Range("K268").Value = 100
$ws.Range("L268").Value = 100
$ws.Range("M268").Value = 10
$ws.Range("N268").Value = 2
$ws.Range("O268").Value = 100
$ws.Range("P268").Value = 100
$ws.Range("Q268").Value = 100
$ws.Range("R268").Value = 100
$ws.Range("S268").Value = 10
$ws.Range("T268").Value = 12
$ws.Range("U268").Value = "Large"
$ws.Range("V268").Value = "Undead"
$ws.Range("W268").Value = "Undead4"
$ws.Range("X268").Value = 1816
$ws.Range("Y268").Value = 240
$ws.Range("Z268").Value = 576
$ws.Range("AA268").Value = 175
$ws.Range("AB268").Value = "Normal"
$ws.Range("AC268").Value = "Elite,Undead"
$ws.Range("AD268").Value = "AiAggressive"
$ws.Range("AE268").Value = 850
$ws.Range("AG268").Value = 0
$ws.Range("AH268").Value = 0.5
$ws.Range("AI268").Value = 1

$ws.Range("AC268").Select()
